# Update forecast output: insert a Week_Start_Date column, fix week labels,
# correct a couple of MyForecast values, switch is_holiday_week to boolean,
# and refresh the dependent Summary totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Forecast Comparison"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column B ("Week_Start_Date"), shifting ASIN .. is_holiday_week
# one column to the right (C..J).
$ws.Columns.Item(2).Insert()

$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Week-start dates for the 16 forecast rows (rows 2-17), written as plain
# text so Excel doesn't coerce them into date serial numbers.
$weekStartDates = @(
    "2025-01-05", "2025-01-12", "2025-01-19", "2025-01-26",
    "2025-02-02", "2025-02-09", "2025-02-16", "2025-02-23",
    "2025-03-02", "2025-03-09", "2025-03-16", "2025-03-23",
    "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20"
)

$ws.Range("B2:B17").NumberFormat = "@"
for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $weekStartDates[$i]
}
$ws.Range("B2:B17").ClearFormats()

# Week labels lose their leading zero (W01 -> W1 ... W9 -> W9; W10.. already
# had no leading zero and stay the same).
$weekLabels = @("W1","W2","W3","W4","W5","W6","W7","W8","W9")
for ($i = 0; $i -lt $weekLabels.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $weekLabels[$i]
}

# Corrected MyForecast values (now column D after the column insert).
$ws.Cells.Item(9, 4).Value = 36
$ws.Cells.Item(16, 4).Value = 34

# is_holiday_week (now column J) becomes a real boolean column.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 10).Value = $false
}

# ---------------------------------------------------------------------------
# Sheet 2: "Summary" - refresh the totals that depend on MyForecast.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B9").Value = "596"
$summary.Range("B10").Value = "312"
$summary.Range("B11").Value = "162"
